# Apply cell updates per the cryptos.xlsx diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new value is a plain string that looks like a number (e.g.
# "224.21") need special handling: assigning such a string straight to
# .Value lets Excel's usual "smart" input parsing turn it into a real
# number (just like it would if a user typed it in), which would change
# the cell type away from the original text cell. To keep it text
# (matching the source workbook, where these are plain string cells, not
# numbers) without touching the cell's style/number-format, we write the
# text via a formula that evaluates to a string, then copy / paste-special
# (values-only) over itself. That bakes the formula result into a plain
# static string cell, with no residual formula and no style change.
function Set-TextValue($range, [string]$text) {
    $escaped = $text.Replace('"', '""')
    $range.Formula = '="' + $escaped + '"'
    $range.Copy()
    $range.PasteSpecial(-4163)
}


$ws.Range('D2').Value = '34.461.08'
$ws.Range('E2').Value = '  -0.13%  '
$ws.Range('D3').Value = '1.801.51'
$ws.Range('E3').Value = '  -0.33%  '
$ws.Range('E4').Value = '  -0.12%  '
Set-TextValue $ws.Range('D5') '224.21'
$ws.Range('E5').Value = '  -1.61%  '
Set-TextValue $ws.Range('D6') '0.601'
$ws.Range('E6').Value = '  +3.57%  '
$ws.Range('E7').Value = '  -0.09%  '
Set-TextValue $ws.Range('D8') '39.26'
$ws.Range('E8').Value = '  +6.63%  '
$ws.Range('E9').Value = '  -3.91%  '
$ws.Range('E10').Value = '  -4.04%  '
Set-TextValue $ws.Range('D11') '0.0983'
$ws.Range('E11').Value = '  +1.93%  '
$ws.Range('D12').Value = '2.061.69'
$ws.Range('E12').Value = '  -0.31%  '
Set-TextValue $ws.Range('D13') '10.90'
$ws.Range('E13').Value = '  -6.11%  '
$ws.Range('D14').Value = '1.803.86'
$ws.Range('E14').Value = '  -0.03%  '
Set-TextValue $ws.Range('D15') '0.629'
$ws.Range('E15').Value = '  -4.04%  '
$ws.Range('D16').Value = '34.432.96'
$ws.Range('E16').Value = '  -0.17%  '
$ws.Range('E17').Value = '  -2.84%  '
Set-TextValue $ws.Range('D18') '67.94'
Set-TextValue $ws.Range('D19') '238.48'
$ws.Range('E19').Value = '  -2.83%  '
$ws.Range('D20').Value = '0.0₃0764'
$ws.Range('E20').Value = '  -3.40%  '
$ws.Range('E21').Value = '  -4.85%  '
$ws.Range('E22').Value = '  -0.13%  '
Set-TextValue $ws.Range('D23') '4.06'
$ws.Range('E23').Value = '  -3.28%  '
$ws.Range('E24').Value = '  -1.62%  '
Set-TextValue $ws.Range('D25') '170.63'
$ws.Range('E25').Value = '  -1.08%  '
Set-TextValue $ws.Range('D26') '17.56'
$ws.Range('E26').Value = '  +3.77%  '
Set-TextValue $ws.Range('D27') '7.63'
$ws.Range('E27').Value = '  -4.49%  '
$ws.Range('E28').Value = '  +2.54%  '
$ws.Range('E29').Value = '  -0.07%  '
$ws.Range('E30').Value = '  -1.80%  '
$ws.Range('E31').Value = '  -2.97%  '
$ws.Range('E32').Value = '  -3.45%  '
Set-TextValue $ws.Range('D33') '3.82'
$ws.Range('E33').Value = '  -5.23%  '
$ws.Range('E34').Value = '  -0.38%  '
$ws.Range('B35').Value = 'TrustWalletToken'
$ws.Range('C35').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
Set-TextValue $ws.Range('D35') '1.05'
$ws.Range('E35').Value = '  -0.60%  '
$ws.Range('B36').Value = 'ImmutableX'
$ws.Range('C36').Value = 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx'
Set-TextValue $ws.Range('D36') '0.638'
$ws.Range('E36').Value = '  -5.29%  '
$ws.Range('D37').Value = '1.300.58'
$ws.Range('E37').Value = '  -6.83%  '
$ws.Range('E38').Value = '  -3.09%  '
$ws.Range('B39').Value = 'HuobiToken'
$ws.Range('C39').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
Set-TextValue $ws.Range('D39') '2.44'
$ws.Range('E39').Value = '  +0.42%  '
$ws.Range('B40').Value = 'RenderToken'
$ws.Range('C40').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue $ws.Range('D40') '2.29'
$ws.Range('E40').Value = '  -6.95%  '
$ws.Range('E41').Value = '  +1.64%  '
Set-TextValue $ws.Range('D42') '81.48'
$ws.Range('E42').Value = '  -2.07%  '
$ws.Range('B43').Value = 'MXToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
Set-TextValue $ws.Range('D43') '2.79'
$ws.Range('E43').Value = '  -1.60%  '
$ws.Range('B44').Value = 'ARBITRUM'
$ws.Range('C44').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
Set-TextValue $ws.Range('D44') '0.944'
$ws.Range('E44').Value = '  -2.80%  '
Set-TextValue $ws.Range('D45') '13.87'
$ws.Range('E45').Value = '  +2.78%  '
Set-TextValue $ws.Range('D46') '0.0516'
$ws.Range('E46').Value = '  +2.98%  '
$ws.Range('D47').Value = '1.963.25'
$ws.Range('E47').Value = '  -0.33%  '
$ws.Range('E49').Value = '  -0.11%  '
Set-TextValue $ws.Range('D50') '101.51'
$ws.Range('E50').Value = '  -2.86%  '
$ws.Range('E51').Value = '  -0.90%  '
